# Crypto list refresh - GitHub Actions scheduled update
# Updates Price (col D) and Volume/1h change (col E) for the rows whose
# market data moved, and re-orders the NEARProtocol / ImmutableX rows
# (35 & 36) to reflect their new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a text value even when it looks like a number
    # (e.g. "1.00", "7.96"), then drop the temporary text format so the
    # cell's style matches the rest of the untouched, unstyled cells.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "62.650.06"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.440.15"
$ws.Range("E3").Value = "  -1.25%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "578.76"
$ws.Range("E5").Value = "  -0.95%  "
Set-TextValue $ws.Range("D6") "147.50"
$ws.Range("E6").Value = "  -0.23%  "
Set-TextValue $ws.Range("D8") "0.481"
$ws.Range("E8").Value = "  +0.37%  "
Set-TextValue $ws.Range("D9") "7.96"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("E10").Value = "  -1.99%  "
Set-TextValue $ws.Range("D11") "0.410"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").Value = "4.023.63"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  +2.37%  "
Set-TextValue $ws.Range("D14") "28.17"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("D15").Value = "3.435.31"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "62.702.09"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  +0.65%  "
Set-TextValue $ws.Range("D19") "14.59"
$ws.Range("E19").Value = "  +1.68%  "
Set-TextValue $ws.Range("D20") "9.03"
$ws.Range("E20").Value = "  -3.37%  "
Set-TextValue $ws.Range("D21") "386.60"
$ws.Range("E21").Value = "  -0.58%  "
Set-TextValue $ws.Range("D22") "75.25"
$ws.Range("E22").Value = "  +0.24%  "
Set-TextValue $ws.Range("D23") "0.560"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "3.582.26"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("E27").Value = "  +0.78%  "
Set-TextValue $ws.Range("D28") "7.63"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  +0.07%  "
Set-TextValue $ws.Range("D30") "7.97"
$ws.Range("E30").Value = "  -3.70%  "
Set-TextValue $ws.Range("D31") "2.11"
$ws.Range("E31").Value = "  -1.67%  "
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -6.49%  "
Set-TextValue $ws.Range("D34") "23.23"
$ws.Range("E34").Value = "  -2.40%  "

# Rows 35 & 36 swap places (NEARProtocol <-> ImmutableX) with refreshed
# price/volume figures for each coin
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "1.62"
$ws.Range("E35").Value = "  +2.60%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "5.32"
$ws.Range("E36").Value = "  -0.23%  "

Set-TextValue $ws.Range("D37") "31.82"
$ws.Range("E37").Value = "  -0.33%  "
Set-TextValue $ws.Range("D38") "6.96"
$ws.Range("E38").Value = "  -2.23%  "
Set-TextValue $ws.Range("D39") "169.97"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "3.469.22"
$ws.Range("E40").Value = "  -1.51%  "
Set-TextValue $ws.Range("D41") "0.0773"
$ws.Range("E41").Value = "  +0.10%  "
Set-TextValue $ws.Range("D42") "0.786"
$ws.Range("E42").Value = "  -2.64%  "
Set-TextValue $ws.Range("D43") "42.56"
$ws.Range("E43").Value = "  +0.54%  "
Set-TextValue $ws.Range("D44") "1.70"
$ws.Range("E44").Value = "  -1.37%  "
Set-TextValue $ws.Range("D45") "4.35"
$ws.Range("E45").Value = "  -2.99%  "
Set-TextValue $ws.Range("D46") "1.18"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "2.563.57"
$ws.Range("E47").Value = "  -2.48%  "
Set-TextValue $ws.Range("D48") "6.95"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  -1.05%  "
Set-TextValue $ws.Range("D50") "22.54"
$ws.Range("E50").Value = "  -3.84%  "
Set-TextValue $ws.Range("D51") "1.00"
$ws.Range("E51").Value = "  +0.03%  "
